# Auto-generated Excel COM-interop edit script
# Applies numeric recalculation updates to the Leve profit tables
# across multiple job sheets (ALC, ARM, CRP, GSM, LTW, WVR), matching
# a scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H12").Value = 200
$ws.Range("I12").Value = 200
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 200
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -30
$ws.Range("N12").ClearContents()

$ws.Range("H62").Value = 5765.3335
$ws.Range("I62").Value = 5940.6
$ws.Range("K62").Value = 5940.6
$ws.Range("M62").Value = -5316.6

$ws.Range("H65").Value = 5765.3335
$ws.Range("I65").Value = 5940.6
$ws.Range("K65").Value = 29703
$ws.Range("M65").Value = -26583

$ws.Range("H70").Value = 1861.1666
$ws.Range("I70").Value = 1572.6666
$ws.Range("J70").Value = 2149.6667
$ws.Range("K70").Value = 4717.9998
$ws.Range("L70").Value = 6449.000100000001
$ws.Range("M70").Value = -4447.9998
$ws.Range("N70").Value = -6989.000100000001

$ws.Range("H73").Value = 1861.1666
$ws.Range("I73").Value = 1572.6666
$ws.Range("J73").Value = 2149.6667
$ws.Range("K73").Value = 4717.9998
$ws.Range("L73").Value = 6449.000100000001
$ws.Range("M73").Value = -3781.9998
$ws.Range("N73").Value = -8321.000100000001

$ws.Range("H76").Value = 6816.3335
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 6816.3335
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 6816.3335
$ws.Range("N76").Value = -7446.3335
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 6816.3335
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 6816.3335
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 6816.3335
$ws.Range("N79").Value = -9000.333500000001
$ws.Range("M79").ClearContents()

$ws.Range("H98").Value = 2407.9565
$ws.Range("I98").Value = 2721.2104
$ws.Range("K98").Value = 2721.2104
$ws.Range("M98").Value = -1223.2104

$ws.Range("H108").Value = 50684
$ws.Range("J108").Value = 50684
$ws.Range("L108").Value = 50684
$ws.Range("N108").Value = -58364

$ws.Range("H122").Value = 2407.9565
$ws.Range("I122").Value = 2721.2104
$ws.Range("K122").Value = 8163.6312
$ws.Range("M122").Value = -5713.6312

$ws.Range("H125").Value = 885.6
$ws.Range("I125").Value = 848
$ws.Range("J125").Value = 1036
$ws.Range("K125").Value = 7632
$ws.Range("L125").Value = 9324
$ws.Range("M125").Value = -5172
$ws.Range("N125").Value = -14244

$ws.Range("H137").Value = 1773.2858
$ws.Range("I137").Value = 1755.8462
$ws.Range("K137").Value = 5267.5386
$ws.Range("M137").Value = -2717.5386


$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H24").Value = 60108.715
$ws.Range("J24").Value = 60108.715
$ws.Range("L24").Value = 60108.715
$ws.Range("N24").Value = -60856.715

$ws.Range("H74").Value = 3147.087
$ws.Range("I74").Value = 2915.4736
$ws.Range("K74").Value = 2915.4736
$ws.Range("M74").Value = -2041.4736

$ws.Range("H77").Value = 3147.087
$ws.Range("I77").Value = 2915.4736
$ws.Range("K77").Value = 14577.368
$ws.Range("M77").Value = -10209.368

$ws.Range("H97").Value = 359.19232
$ws.Range("I97").Value = 340.5
$ws.Range("J97").Value = 462
$ws.Range("K97").Value = 340.5
$ws.Range("L97").Value = 462
$ws.Range("M97").Value = 155.5
$ws.Range("N97").Value = -1454

$ws.Range("H100").Value = 60108.715
$ws.Range("J100").Value = 60108.715
$ws.Range("L100").Value = 60108.715
$ws.Range("N100").Value = -62272.715

$ws.Range("H102").Value = 16146444
$ws.Range("I102").Value = 20001510
$ws.Range("K102").Value = 20001510
$ws.Range("M102").Value = -19999888

$ws.Range("H122").Value = 6902909
$ws.Range("I122").Value = 7698683
$ws.Range("K122").Value = 23096049
$ws.Range("M122").Value = -23093599

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()


$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H7").Value = 89.15000000000001
$ws.Range("I7").Value = 77.42856999999999
$ws.Range("K7").Value = 77.42856999999999
$ws.Range("M7").Value = 35.57143000000001

$ws.Range("H58").Value = 4377.227
$ws.Range("I58").Value = 4589.25
$ws.Range("J58").Value = 2257
$ws.Range("K58").Value = 4589.25
$ws.Range("L58").Value = 2257
$ws.Range("M58").Value = -4386.25
$ws.Range("N58").Value = -2663

$ws.Range("H132").Value = 2786.75
$ws.Range("I132").Value = 2131.0908
$ws.Range("K132").Value = 6393.2724
$ws.Range("M132").Value = -3863.2724

$ws.Range("H136").Value = 4377.227
$ws.Range("I136").Value = 4589.25
$ws.Range("J136").Value = 2257
$ws.Range("K136").Value = 13767.75
$ws.Range("L136").Value = 6771
$ws.Range("M136").Value = -11217.75
$ws.Range("N136").Value = -11871


$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H102").Value = 1298.1111
$ws.Range("I102").Value = 1298.1111
$ws.Range("K102").Value = 1298.1111
$ws.Range("M102").Value = 323.8888999999999

$ws.Range("H113").Value = 3833.3333
$ws.Range("J113").Value = 4850
$ws.Range("L113").Value = 4850
$ws.Range("N113").Value = -9190

$ws.Range("H122").Value = 35715730
$ws.Range("I122").Value = 43479388
$ws.Range("J122").Value = 2880.8
$ws.Range("K122").Value = 130438164
$ws.Range("L122").Value = 8642.400000000001
$ws.Range("M122").Value = -130435714
$ws.Range("N122").Value = -13542.4

$ws.Range("H126").Value = 6677.4
$ws.Range("I126").Value = 6677.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 20032.2
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -17562.2
$ws.Range("N126").ClearContents()


$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 5481.615
$ws.Range("J7").Value = 7399
$ws.Range("L7").Value = 7399
$ws.Range("N7").Value = -7623

$ws.Range("H22").Value = 4999.5
$ws.Range("J22").Value = 5499.6665
$ws.Range("L22").Value = 5499.6665
$ws.Range("N22").Value = -6089.6665

$ws.Range("H27").Value = 4999.5
$ws.Range("J27").Value = 5499.6665
$ws.Range("L27").Value = 5499.6665
$ws.Range("N27").Value = -5713.6665

$ws.Range("H40").Value = 7872.9287
$ws.Range("J40").Value = 6723.857
$ws.Range("L40").Value = 6723.857
$ws.Range("N40").Value = -6995.857

$ws.Range("H68").Value = 6300.4287
$ws.Range("I68").Value = 1987
$ws.Range("K68").Value = 1987
$ws.Range("M68").Value = -1238

$ws.Range("H71").Value = 6300.4287
$ws.Range("I71").Value = 1987
$ws.Range("K71").Value = 9935
$ws.Range("M71").Value = -6191

$ws.Range("H101").Value = 75876.836
$ws.Range("J101").Value = 75876.836
$ws.Range("L101").Value = 75876.836
$ws.Range("N101").Value = -82366.836

$ws.Range("H122").Value = 8008.4473
$ws.Range("I122").Value = 7904.9688
$ws.Range("K122").Value = 23714.9064
$ws.Range("M122").Value = -21264.9064

$ws.Range("H126").Value = 5481.615
$ws.Range("J126").Value = 7399
$ws.Range("L126").Value = 22197
$ws.Range("N126").Value = -27137

$ws.Range("H136").Value = 5461.3823
$ws.Range("I136").Value = 4421.1035
$ws.Range("K136").Value = 13263.3105
$ws.Range("M136").Value = -10713.3105


$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H81").Value = 4201.4443
$ws.Range("I81").Value = 1433
$ws.Range("K81").Value = 2866
$ws.Range("M81").Value = -1805

$ws.Range("H84").Value = 4201.4443
$ws.Range("I84").Value = 1433
$ws.Range("K84").Value = 14330
$ws.Range("M84").Value = -9026

$ws.Range("H122").Value = 3284.4
$ws.Range("I122").Value = 3301.4
$ws.Range("J122").Value = 3250.4
$ws.Range("K122").Value = 9904.200000000001
$ws.Range("L122").Value = 9751.200000000001
$ws.Range("M122").Value = -7454.200000000001
$ws.Range("N122").Value = -14651.2

$ws.Range("H126").Value = 2331.2
$ws.Range("I126").Value = 2194.125
$ws.Range("K126").Value = 6582.375
$ws.Range("M126").Value = -4112.375

$ws.Range("H132").Value = 1913.5
$ws.Range("I132").Value = 1884.8823
$ws.Range("K132").Value = 5654.6469
$ws.Range("M132").Value = -3124.6469

$ws.Range("H136").Value = 4461.2188
$ws.Range("I136").Value = 3836.6538
$ws.Range("K136").Value = 11509.9614
$ws.Range("M136").Value = -8959.9614

